# Manage DB + sync fix
# - Rename a few headers, update existing ID values, append two new result
#   rows (Meow / Mona), drop the header-row border, and resize a few columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels -------------------------------------------------
$ws.Range("D1").Value = "Score (%)"
$ws.Range("G1").Value = "Time (sec)"
$ws.Range("H1").Value = "Completed At"

# --- Renumber existing ID column values ------------------------------------
$ws.Range("A2").Value = 105
$ws.Range("A3").Value = 106
$ws.Range("A4").Value = 107
$ws.Range("A5").Value = 108

# --- Append new rows of quiz results ---------------------------------------
$ws.Cells.Item(6, 1).Value = 109
$ws.Cells.Item(6, 2).Value = "Meow"
$ws.Cells.Item(6, 3).Value = "I07547"
$ws.Cells.Item(6, 4).Value = 40
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Value = 12
$ws.Cells.Item(6, 8).Value = "2025-04-27 21:41:00"

$ws.Cells.Item(7, 1).Value = 110
$ws.Cells.Item(7, 2).Value = "Mona"
$ws.Cells.Item(7, 3).Value = "I332232"
$ws.Cells.Item(7, 4).Value = 40
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(7, 7).Value = 23
$ws.Cells.Item(7, 8).Value = "2025-04-27 21:42:00"

# --- Drop the border that used to sit under the header row -----------------
$ws.Range("A1:H1").Borders.LineStyle = 0

# --- Resize columns to the new "best fit" widths ----------------------------
$ws.Columns.Item(1).ColumnWidth = 3.6666666666666665
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666
